$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is purely numeric-looking text (e.g. "600.62") need
# NumberFormat forced to text ("@") first, otherwise Excel auto-converts the
# assignment into a real number (losing the original text-cell semantics and
# introducing float rounding). Values with two dots (e.g. "65.140.00") or a
# percent sign are never parsed as numbers by Excel, so they are safe as-is.

$ws.Range("D2").Value = '65.140.00'
$ws.Range("E2").Value = '  -0.87%  '
$ws.Range("D3").Value = '3.531.85'
$ws.Range("E3").Value = '  +2.43%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '600.62'
$ws.Range("E5").Value = '  +1.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.85'
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("D7").Value = '3.534.91'
$ws.Range("E7").Value = '  +2.62%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  -2.93%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.124'
$ws.Range("E10").Value = '  +1.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.92'
$ws.Range("E11").Value = '  -5.99%  '
$ws.Range("E12").Value = '  +2.48%  '
$ws.Range("D13").Value = '4.133.89'
$ws.Range("E13").Value = '  +2.65%  '
$ws.Range("E14").Value = '  +1.19%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.11'
$ws.Range("E15").Value = '  +1.76%  '
$ws.Range("D16").Value = '3.539.39'
$ws.Range("E16").Value = '  +2.78%  '
$ws.Range("E17").Value = '  +1.50%  '
$ws.Range("D18").Value = '65.263.65'
$ws.Range("E18").Value = '  -0.53%  '
$ws.Range("E19").Value = '  +4.29%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.94'
$ws.Range("E20").Value = '  +0.23%  '
$ws.Range("E21").Value = '  +3.52%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '392.35'
$ws.Range("E22").Value = '  -0.78%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.572'
$ws.Range("E23").Value = '  +2.87%  '
$ws.Range("D24").Value = '3.673.63'
$ws.Range("E24").Value = '  +2.37%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '73.76'
$ws.Range("E25").Value = '  +0.31%  '
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("E27").Value = '  +6.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.69'
$ws.Range("E28").Value = '  +6.45%  '
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("E30").Value = '  +1.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.20'
$ws.Range("E31").Value = '  -1.72%  '
$ws.Range("D32").Value = '3.545.44'
$ws.Range("E32").Value = '  +2.73%  '
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.73'
$ws.Range("E34").Value = '  +2.75%  '
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("E36").Value = '  +5.73%  '
$ws.Range("E37").Value = '  -0.38%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '168.54'
$ws.Range("E38").Value = '  -2.59%  '
$ws.Range("E39").Value = '  +3.70%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.96'
$ws.Range("E40").Value = '  +2.59%  '
$ws.Range("E41").Value = '  +4.31%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.823'
$ws.Range("E42").Value = '  -0.41%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.33'
$ws.Range("E43").Value = '  +12.64%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.77'
$ws.Range("E44").Value = '  -2.49%  '
$ws.Range("E45").Value = '  +0.15%  '
$ws.Range("E46").Value = '  -0.64%  '
$ws.Range("E47").Value = '  +2.20%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.18'
$ws.Range("E48").Value = '  +4.73%  '
$ws.Range("D49").Value = '2.412.68'
$ws.Range("E49").Value = '  +8.88%  '
$ws.Range("E50").Value = '  +2.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '301.42'
$ws.Range("E51").Value = '  +6.27%  '
